$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2030651340996169
$ws.Range("C2").Value = 0.5747126436781609
$ws.Range("J2").Value = 0.007662835249042145
$ws.Range("P2").Value = 0.1532567049808429
$ws.Range("S2").Value = 0.06130268199233716
$ws.Range("B3").Value = 0.02547770700636943
$ws.Range("C3").Value = 0.03821656050955414
$ws.Range("J3").Value = 0.03821656050955414
$ws.Range("P3").Value = 0.7388535031847133
$ws.Range("S3").Value = 0.1592356687898089
$ws.Range("J4").Value = 0.05714285714285714
$ws.Range("P4").Value = 0.6571428571428571
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.08333333333333333
$ws.Range("D6").Value = 0.01041666666666667
$ws.Range("J6").Value = 0.2552083333333333
$ws.Range("O6").Value = 0.015625
$ws.Range("Q6").Value = 0.1145833333333333
$ws.Range("R6").Value = 0.08333333333333333
$ws.Range("S6").Value = 0.375
$ws.Range("B7").Value = 0.1369047619047619
$ws.Range("D7").Value = 0.01785714285714286
$ws.Range("F7").Value = 0.05357142857142857
$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("O7").Value = 0.005952380952380952
$ws.Range("Q7").Value = 0.1428571428571428
$ws.Range("R7").Value = 0.07142857142857142
$ws.Range("S7").Value = 0.4285714285714285
$ws.Range("B8").Value = 0.1195928753180662
$ws.Range("D8").Value = 0.01272264631043257
$ws.Range("E8").Value = 0.002544529262086514
$ws.Range("F8").Value = 0.04325699745547074
$ws.Range("J8").Value = 0.1094147582697201
$ws.Range("O8").Value = 0.01272264631043257
$ws.Range("Q8").Value = 0.2061068702290076
$ws.Range("R8").Value = 0.08651399491094147
$ws.Range("S8").Value = 0.4071246819338422
$ws.Range("B9").Value = 0.1304347826086956
$ws.Range("D9").Value = 0.01739130434782609
$ws.Range("E9").Value = 0.008695652173913044
$ws.Range("F9").Value = 0.06956521739130435
$ws.Range("J9").Value = 0.1478260869565217
$ws.Range("Q9").Value = 0.1565217391304348
$ws.Range("R9").Value = 0.06956521739130435
$ws.Range("S9").Value = 0.4
$ws.Range("B10").Value = 0.09208103130755065
$ws.Range("D10").Value = 0.02394106813996317
$ws.Range("E10").Value = 0.0009208103130755065
$ws.Range("F10").Value = 0.0718232044198895
$ws.Range("J10").Value = 0.1243093922651934
$ws.Range("O10").Value = 0.01197053406998158
$ws.Range("Q10").Value = 0.2274401473296501
$ws.Range("R10").Value = 0.07642725598526703
$ws.Range("S10").Value = 0.3710865561694291
$ws.Range("G11").Value = 0.15
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.2230769230769231
$ws.Range("L11").Value = 0.5230769230769231
$ws.Range("S11").Value = 0.003846153846153846
$ws.Range("G12").Value = 0.7681159420289855
$ws.Range("J12").Value = 0.1956521739130435
$ws.Range("K12").Value = 0.007246376811594203
$ws.Range("L12").Value = 0.01449275362318841
$ws.Range("S12").Value = 0.01449275362318841
$ws.Range("G13").Value = 0.6904761904761905
$ws.Range("J13").Value = 0.3095238095238095
$ws.Range("F15").Value = 0.02958579881656805
$ws.Range("H15").Value = 0.1479289940828402
$ws.Range("I15").Value = 0.04733727810650887
$ws.Range("J15").Value = 0.3846153846153846
$ws.Range("K15").Value = 0.1183431952662722
$ws.Range("M15").Value = 0.01775147928994083
$ws.Range("O15").Value = 0.0650887573964497
$ws.Range("S15").Value = 0.1893491124260355
$ws.Range("F16").Value = 0.03932584269662921
$ws.Range("H16").Value = 0.151685393258427
$ws.Range("I16").Value = 0.07303370786516854
$ws.Range("J16").Value = 0.4550561797752809
$ws.Range("K16").Value = 0.0898876404494382
$ws.Range("M16").Value = 0.01685393258426966
$ws.Range("O16").Value = 0.05617977528089887
$ws.Range("S16").Value = 0.1179775280898876
$ws.Range("F17").Value = 0.0310880829015544
$ws.Range("H17").Value = 0.1580310880829016
$ws.Range("I17").Value = 0.07772020725388601
$ws.Range("J17").Value = 0.4404145077720207
$ws.Range("K17").Value = 0.1088082901554404
$ws.Range("M17").Value = 0.01813471502590673
$ws.Range("O17").Value = 0.04404145077720207
$ws.Range("S17").Value = 0.1217616580310881
$ws.Range("F18").Value = 0.006578947368421052
$ws.Range("H18").Value = 0.2171052631578947
$ws.Range("I18").Value = 0.03289473684210526
$ws.Range("J18").Value = 0.4276315789473684
$ws.Range("K18").Value = 0.05263157894736842
$ws.Range("M18").Value = 0.03289473684210526
$ws.Range("O18").Value = 0.07894736842105263
$ws.Range("S18").Value = 0.1513157894736842
$ws.Range("F19").Value = 0.01937984496124031
$ws.Range("H19").Value = 0.2354651162790698
$ws.Range("I19").Value = 0.0562015503875969
$ws.Range("J19").Value = 0.3662790697674418
$ws.Range("K19").Value = 0.1114341085271318
$ws.Range("M19").Value = 0.0251937984496124
$ws.Range("O19").Value = 0.07170542635658915
$ws.Range("S19").Value = 0.1143410852713178
